$d = $word.ActiveDocument

# Helper: split the run at an absolute document character offset by
# briefly adding then deleting a bookmark there. Word (and this engine)
# splits the underlying run at a bookmark boundary and leaves it split
# even after the bookmark itself is removed, while runs keep identical
# formatting to their neighbour.
function Split-RunAt($pos) {
    $mark = $d.Range($pos, $pos)
    $d.Bookmarks.Add("tmp_split_pt", $mark) | Out-Null
    $d.Bookmarks("tmp_split_pt").Delete()
}

# ---------------------------------------------------------------------
# 1. Job title line: "NETWORK SUPPORT SPECIALIST – LEAD, " becomes four
#    runs reading "SR." " LEAD" " - TECHNOLOGY," " " (promotion title).
# ---------------------------------------------------------------------
$titleRng = $d.Content
$titleFound = $titleRng.Find.Execute("NETWORK SUPPORT SPECIALIST " + [char]0x2013 + " LEAD, ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($titleFound) {
    $titleStart = $titleRng.Start
    $titleRng.Text = "SR. LEAD - TECHNOLOGY, "

    Split-RunAt ($titleStart + 3)
    Split-RunAt ($titleStart + 8)
    Split-RunAt ($titleStart + 22)
}

# ---------------------------------------------------------------------
# 2. First bullet under that role.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Lead network troubleshooting for smart home IoT systems", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Collaborate across development teams to streamline ticket management", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Second bullet under that role, which also relocates the `_GoBack`
#    bookmark so it sits between " in-house" and " technicians".
# ---------------------------------------------------------------------
$bulletRng = $d.Content
$bulletFound = $bulletRng.Find.Execute("Train technicians regarding wireless networking principles", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($bulletFound) {
    $bulletStart = $bulletRng.Start
    $newBulletText = "Analyze data produced by in-house technicians to drive efficiency and inform best practice"
    $bulletRng.Text = $newBulletText

    $split1 = $bulletStart + "Analyze data produced by".Length
    $split2 = $bulletStart + "Analyze data produced by in-house".Length

    Split-RunAt $split1
    Split-RunAt $split2

    # Remove the bookmark from its old home in the Education section
    # (it sat right before "GPA: 3.46") and re-create it at the new spot.
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $newMark = $d.Range($split2, $split2)
    $d.Bookmarks.Add("_GoBack", $newMark) | Out-Null
}
